$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "71.254.79"
$ws.Range("E2").Value = "  +6.58%  "
$ws.Range("D3").Value = "3.776.28"
$ws.Range("E3").Value = "  +22.54%  "
$ws.Range("E4").Value = "  -0.13%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "613.22"
$ws.Range("E5").Value = "  +7.64%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "179.77"
$ws.Range("E6").Value = "  +1.69%  "
$ws.Range("D7").Value = "3.775.76"
$ws.Range("E7").Value = "  +22.51%  "
$ws.Range("E8").Value = "  -0.08%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.542"
$ws.Range("E9").Value = "  +6.09%  "
$ws.Range("E10").Value = "  +10.41%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.45"
$ws.Range("E11").Value = "  +1.50%  "
$ws.Range("E12").Value = "  +7.86%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "40.53"
$ws.Range("E13").Value = "  +13.10%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.0000257"
$ws.Range("E14").Value = "  +7.83%  "
$ws.Range("D15").Value = "4.396.27"
$ws.Range("E15").Value = "  +22.09%  "
$ws.Range("D16").Value = "3.772.51"
$ws.Range("E16").Value = "  +22.23%  "
$ws.Range("D17").Value = "71.314.59"
$ws.Range("E17").Value = "  +6.67%  "
$ws.Range("E18").Value = "  +1.39%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.53"
$ws.Range("E19").Value = "  +7.95%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "524.23"
$ws.Range("E20").Value = "  +8.59%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "16.72"
$ws.Range("E21").Value = "  +1.62%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "9.35"
$ws.Range("E22").Value = "  +22.14%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.746"
$ws.Range("E23").Value = "  +9.51%  "
$ws.Range("B24").Value = "Litecoin"
$ws.Range("C24").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "88.46"
$ws.Range("E24").Value = "  +6.24%  "
$ws.Range("B25").Value = "Fetch.AI"
$ws.Range("C25").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.49"
$ws.Range("E25").Value = "  +11.52%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "13.53"
$ws.Range("E26").Value = "  +8.02%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "11.04"
$ws.Range("E27").Value = "  +9.49%  "
$ws.Range("E28").Value = "  +0.03%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.53"
$ws.Range("E29").Value = "  +11.05%  "
$ws.Range("B30").Value = "PEPE"
$ws.Range("C30").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0000116"
$ws.Range("E30").Value = "  +24.83%  "
$ws.Range("B31").Value = "PancakeSwap"
$ws.Range("C31").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.91"
$ws.Range("E31").Value = "  +13.01%  "
$ws.Range("B32").Value = "NEARProtocol"
$ws.Range("C32").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "8.05"
$ws.Range("E32").Value = "  +2.92%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "32.09"
$ws.Range("E33").Value = "  +15.02%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.115"
$ws.Range("E34").Value = "  +3.94%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.00"
$ws.Range("E35").Value = "  -0.01%  "
$ws.Range("E36").Value = "  +12.73%  "
$ws.Range("E37").Value = "  +10.73%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.23"
$ws.Range("E38").Value = "  +11.82%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.341"
$ws.Range("E39").Value = "  +10.44%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.133"
$ws.Range("E40").Value = "  +7.81%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "51.89"
$ws.Range("E41").Value = "  +6.03%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "432.34"
$ws.Range("E42").Value = "  +18.23%  "
$ws.Range("D43").Value = "3.158.34"
$ws.Range("E43").Value = "  +13.07%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.82"
$ws.Range("E44").Value = "  +7.58%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "43.85"
$ws.Range("E45").Value = "  -6.70%  "
$ws.Range("E46").Value = "  +3.77%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0367"
$ws.Range("E47").Value = "  +7.25%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "27.90"
$ws.Range("E48").Value = "  +10.31%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "140.73"
$ws.Range("E49").Value = "  +3.89%  "
$ws.Range("E50").Value = "  +0.01%  "
